$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.276.89"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.61"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.59"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4621"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.86"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9753"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.70"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.27"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.793"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.081"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07060"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.74"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009742"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.02"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.273.42"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.468"
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.159.51"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.091"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.39"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.772"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.95"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.844"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09323"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8645"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.160"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.300"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05777"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.156"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.639"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5656"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000003078"
$ws.Range("E42").Value = "  +11.83%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1781"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.396"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.711"
$ws.Range("E45").Value = "  +7.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5262"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.45"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06858"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.813"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.12"
$ws.Range("E51").Value = "  -1.71%  "
